# Auto-generated Excel COM-interop script
# Applies numeric corrections to market-price / profit columns (H-N)
# across several Leve-profit worksheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000.1111
$ws.Range("I18").Value = 1000.1111
$ws.Range("K18").Value = 1000.1111
$ws.Range("M18").Value = -716.1111
$ws.Range("H21").Value = 49211
$ws.Range("I21").Value = 41513.75
$ws.Range("J21").Value = 80000
$ws.Range("K21").Value = 41513.75
$ws.Range("L21").Value = 80000
$ws.Range("M21").Value = -41045.75
$ws.Range("N21").Value = -80936
$ws.Range("H23").Value = 49211
$ws.Range("I23").Value = 41513.75
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 41513.75
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -41279.75
$ws.Range("N23").Value = -80468
$ws.Range("H34").Value = 2272
$ws.Range("I34").Value = 2272
$ws.Range("K34").Value = 2272
$ws.Range("M34").Value = -2069
$ws.Range("H36").Value = 2272
$ws.Range("I36").Value = 2272
$ws.Range("K36").Value = 2272
$ws.Range("M36").Value = -1557
$ws.Range("H125").Value = 2948.1667
$ws.Range("J125").Value = 3264.2222
$ws.Range("L125").Value = 29377.9998
$ws.Range("N125").Value = -34297.99980000001
$ws.Range("H137").Value = 19126746
$ws.Range("I137").Value = 3402318.8
$ws.Range("J137").Value = 83334824
$ws.Range("K137").Value = 10206956.4
$ws.Range("L137").Value = 250004472
$ws.Range("M137").Value = -10204406.4
$ws.Range("N137").Value = -250009572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10419409
$ws.Range("I61").Value = 12823463
$ws.Range("J61").Value = 1839.6666
$ws.Range("K61").Value = 12823463
$ws.Range("L61").Value = 1839.6666
$ws.Range("M61").Value = -12823251
$ws.Range("N61").Value = -2263.6666
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 3281545.5
$ws.Range("I74").Value = 4348314.5
$ws.Range("J74").Value = 10120.066
$ws.Range("K74").Value = 4348314.5
$ws.Range("L74").Value = 10120.066
$ws.Range("M74").Value = -4347440.5
$ws.Range("N74").Value = -11868.066
$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35676
$ws.Range("H77").Value = 3281545.5
$ws.Range("I77").Value = 4348314.5
$ws.Range("J77").Value = 10120.066
$ws.Range("K77").Value = 21741572.5
$ws.Range("L77").Value = 50600.33
$ws.Range("M77").Value = -21737204.5
$ws.Range("N77").Value = -59336.33
$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37340
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 1549.9166
$ws.Range("I122").Value = 1529.2632
$ws.Range("J122").Value = 1628.4
$ws.Range("K122").Value = 4587.7896
$ws.Range("L122").Value = 4885.200000000001
$ws.Range("M122").Value = -2137.7896
$ws.Range("N122").Value = -9785.200000000001
$ws.Range("H132").Value = 635179.9399999999
$ws.Range("I132").Value = 719879
$ws.Range("J132").Value = 126985.375
$ws.Range("K132").Value = 2159637
$ws.Range("L132").Value = 380956.125
$ws.Range("M132").Value = -2157107
$ws.Range("N132").Value = -386016.125
$ws.Range("H136").Value = 10419409
$ws.Range("I136").Value = 12823463
$ws.Range("J136").Value = 1839.6666
$ws.Range("K136").Value = 38470389
$ws.Range("L136").Value = 5518.9998
$ws.Range("M136").Value = -38467839
$ws.Range("N136").Value = -10618.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 15737272
$ws.Range("I6").Value = 18360068
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 18360068
$ws.Range("L6").Value = 500
$ws.Range("M6").Value = -18359955
$ws.Range("N6").Value = -726
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H43").Value = 30385.334
$ws.Range("J43").Value = 30385.334
$ws.Range("L43").Value = 30385.334
$ws.Range("N43").Value = -30753.334
$ws.Range("H58").Value = 1217.625
$ws.Range("I58").Value = 1021.6667
$ws.Range("J58").Value = 2141.4285
$ws.Range("K58").Value = 1021.6667
$ws.Range("L58").Value = 2141.4285
$ws.Range("M58").Value = -818.6667
$ws.Range("N58").Value = -2547.4285
$ws.Range("H101").Value = 30385.334
$ws.Range("J101").Value = 30385.334
$ws.Range("L101").Value = 30385.334
$ws.Range("N101").Value = -36875.334
$ws.Range("H134").Value = 1850.6154
$ws.Range("I134").Value = 1912
$ws.Range("K134").Value = 5736
$ws.Range("M134").Value = -3201
$ws.Range("H136").Value = 1217.625
$ws.Range("I136").Value = 1021.6667
$ws.Range("J136").Value = 2141.4285
$ws.Range("K136").Value = 3065.0001
$ws.Range("L136").Value = 6424.2855
$ws.Range("M136").Value = -515.0001000000002
$ws.Range("N136").Value = -11524.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11112120
$ws.Range("I5").Value = 1477.7
$ws.Range("J5").Value = 14286589
$ws.Range("K5").Value = 4433.1
$ws.Range("L5").Value = 42859767
$ws.Range("M5").Value = -4321.1
$ws.Range("N5").Value = -42859991
$ws.Range("H68").Value = 1885.9697
$ws.Range("I68").Value = 896.1667
$ws.Range("J68").Value = 2451.5715
$ws.Range("K68").Value = 2688.5001
$ws.Range("L68").Value = 7354.7145
$ws.Range("M68").Value = -1877.5001
$ws.Range("N68").Value = -8976.7145
$ws.Range("H71").Value = 1885.9697
$ws.Range("I71").Value = 896.1667
$ws.Range("J71").Value = 2451.5715
$ws.Range("K71").Value = 8065.5003
$ws.Range("L71").Value = 22064.1435
$ws.Range("M71").Value = -4009.5003
$ws.Range("N71").Value = -30176.1435
$ws.Range("H106").Value = 5483.1816
$ws.Range("J106").Value = 5483.1816
$ws.Range("L106").Value = 16449.5448
$ws.Range("N106").Value = -18341.5448
$ws.Range("H122").Value = 22322954
$ws.Range("I122").Value = 30303832
$ws.Range("J122").Value = 4765025
$ws.Range("K122").Value = 272734488
$ws.Range("L122").Value = 42885225
$ws.Range("M122").Value = -272732038
$ws.Range("N122").Value = -42890125
$ws.Range("H135").Value = 11112120
$ws.Range("I135").Value = 1477.7
$ws.Range("J135").Value = 14286589
$ws.Range("K135").Value = 13299.3
$ws.Range("L135").Value = 128579301
$ws.Range("M135").Value = -10764.3
$ws.Range("N135").Value = -128584371

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 49325.2
$ws.Range("J51").Value = 49325.2
$ws.Range("L51").Value = 49325.2
$ws.Range("N51").Value = -50343.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1716.6666
$ws.Range("I122").Value = 1775
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 5325
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -2875
$ws.Range("N122").Value = -9700
$ws.Range("H132").Value = 82044.42999999999
$ws.Range("I132").Value = 94718.664
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 284155.992
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -281625.992
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 14301
$ws.Range("J101").Value = 14301
$ws.Range("L101").Value = 14301
$ws.Range("N101").Value = -20791
$ws.Range("H132").Value = 4122.8223
$ws.Range("I132").Value = 4277.5366
$ws.Range("J132").Value = 2537
$ws.Range("K132").Value = 12832.6098
$ws.Range("L132").Value = 7611
$ws.Range("M132").Value = -10302.6098
$ws.Range("N132").Value = -12671

